$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 13797
$ws1.Range("F4").Value = 39
$ws1.Range("F5").Value = 79
$ws1.Range("F7").Value = 2197
$ws1.Range("F8").Value = 204
$ws1.Range("F9").Value = 136
$ws1.Range("F10").Value = 118
$ws1.Range("F11").Value = 250
$ws1.Range("F13").Value = 616
$ws1.Range("F14").Value = 467
$ws1.Range("F15").Value = 533
$ws1.Range("F16").Value = 349
$ws1.Range("F18").Value = 318
$ws1.Range("F19").Value = 903
$ws1.Range("F20").Value = 167
$ws1.Range("F21").Value = 93
$ws1.Range("F22").Value = 44
$ws1.Range("F23").Value = 2
$ws1.Range("F24").Value = 11
$ws1.Range("F25").Value = 115
$ws1.Range("F26").Value = 46

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 59
$ws2.Range("F4").Value = 144
$ws2.Range("F5").Value = 75
$ws2.Range("F6").Value = 134
$ws2.Range("F7").Value = 186
$ws2.Range("F8").Value = 2117
$ws2.Range("F13").Value = 89
$ws2.Range("F15").Value = 1904

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 238
$ws3.Range("F3").Value = 198
$ws3.Range("F4").Value = 139

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 238
$ws4.Range("F3").Value = 13797
$ws4.Range("F5").Value = 39
$ws4.Range("F6").Value = 79
$ws4.Range("F8").Value = 59
$ws4.Range("F10").Value = 2197
$ws4.Range("F11").Value = 198
$ws4.Range("F12").Value = 204
$ws4.Range("F13").Value = 136
$ws4.Range("F14").Value = 118
$ws4.Range("F15").Value = 250
$ws4.Range("F16").Value = 144
$ws4.Range("F18").Value = 75
$ws4.Range("F19").Value = 134
$ws4.Range("F20").Value = 139
$ws4.Range("F21").Value = 616
$ws4.Range("F22").Value = 467
$ws4.Range("F23").Value = 533
$ws4.Range("F24").Value = 349
$ws4.Range("F26").Value = 318
$ws4.Range("F27").Value = 903
$ws4.Range("F28").Value = 186
$ws4.Range("F29").Value = 2117
$ws4.Range("F34").Value = 167
$ws4.Range("F35").Value = 93
$ws4.Range("F36").Value = 44
$ws4.Range("F37").Value = 2
$ws4.Range("F38").Value = 89
$ws4.Range("F40").Value = 11
$ws4.Range("F41").Value = 115
$ws4.Range("F42").Value = 46
$ws4.Range("F43").Value = 1904
